# Migrate mountebank stub example to use api template framework.
# Add flag for mountebank api testing.
# Add create and delete stub functions in Before and After annotations.
#
# This adds a new "testMountebankAPI" Before/After block to the API sheet
# (header row, one data row, footer row), matching the existing
# testAPIWithXML-style blocks already on the sheet, and updates the
# active sheet / selection / zoom to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")
$jdbc = $wb.Worksheets.Item("JDBC")

# --- Copy cell formatting for the new rows from the most similar existing
#     block (row 12 = header row style, row 13 = data row style) so the
#     new cells pick up the same style index (s="8") as their siblings. ---
$ws.Range("B12:F12").Copy()
$ws.Range("B16:F16").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B13:E13").Copy()
$ws.Range("B17:E17").PasteSpecial(-4122)   # xlPasteFormats

# --- Thin blank separator row above the new block (matches the blank
#     rows that precede every other Before/After block on this sheet). ---
$ws.Range("A15").EntireRow.RowHeight = 12

# --- Values: set in the same order they were first typed so new shared
#     strings land in the expected order ("test", "4", "testMountebankAPI"). ---
$ws.Cells.Item(17, 4).Value = "test"               # D17
$ws.Cells.Item(17, 2).Value = "4"                  # B17
$ws.Cells.Item(16, 1).Value = "testMountebankAPI"  # A16 (Before marker / test name)
$ws.Cells.Item(18, 7).Value = "testMountebankAPI"  # G18 (After marker / test name)

# Row 16 - header
$ws.Cells.Item(16, 2).Value = "TemplateId"
$ws.Cells.Item(16, 3).Value = "type"
$ws.Cells.Item(16, 4).Value = "jsonPath"
$ws.Cells.Item(16, 5).Value = "requestType"
$ws.Cells.Item(16, 6).Value = "response"

# Row 17 - data
$ws.Cells.Item(17, 3).Value = "JSON"
$ws.Cells.Item(17, 5).Value = "POST"
$ws.Cells.Item(17, 6).Value = 400

# --- View state: author ended up on the API sheet, zoomed in, with G18
#     selected; JDBC (previously active) is zoomed out with B8 selected. ---
$jdbc.Activate()
$excel.Windows.Item(1).Zoom = 200
$jdbc.Range("B8").Select()

$ws.Activate()
$excel.Windows.Item(1).Zoom = 248
$ws.Range("G18").Select()
